$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.678.70'
$ws.Range('E2').Value = '  -3.85%  '
$ws.Range('D3').Value = '3.311.46'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.33'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '182.37'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.73%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.599'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.31%  '
$ws.Range('E9').Value = '  -3.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.63'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.62%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.404'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.20%  '
$ws.Range('D12').Value = '3.887.33'
$ws.Range('E12').Value = '  -0.96%  '
$ws.Range('E13').Value = '  -0.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.14'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.75%  '
$ws.Range('D15').Value = '66.734.32'
$ws.Range('E15').Value = '  -3.73%  '
$ws.Range('D17').Value = '3.320.41'
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('E18').Value = '  -0.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '434.24'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('E20').Value = '  -2.44%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.63'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.79'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.90%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('E24').Value = '  -0.40%  '
$ws.Range('E25').Value = '  -2.55%  '
$ws.Range('E26').Value = '  +1.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.05'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.34%  '
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('E29').Value = '  -2.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.80'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.33'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.61%  '
$ws.Range('E33').Value = '  -3.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.23'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.84%  '
$ws.Range('E35').Value = '  -1.36%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '160.28'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.66%  '
$ws.Range('E37').Value = '  -3.76%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '27.28'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.81%  '
$ws.Range('D39').Value = '2.814.66'
$ws.Range('E39').Value = '  +2.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.791'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.23%  '
$ws.Range('E41').Value = '  -2.62%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.25'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.80%  '
$ws.Range('E43').Value = '  -1.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.17'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '24.36'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.49%  '
$ws.Range('E46').Value = '  -6.95%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '320.25'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.90%  '
$ws.Range('E48').Value = '  -3.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.984'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.27%  '
$ws.Range('E50').Value = '  -1.87%  '
$ws.Range('E51').Value = '  -1.69%  '
